$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header cells
$ws.Range("C1").Value = "rules"
$ws.Range("D1").Value = "adaptive_filter"

# Row 2
$ws.Range("D2").Value = "wRLS"
$ws.Range("E2").Value = 144.8447666333835
$ws.Range("F2").Value = 0.3070031139977897
$ws.Range("G2").Value = 107.9718208199604

# Row 3
$ws.Range("D3").Value = "wRLS"
$ws.Range("E3").Value = 143.7238371389704
$ws.Range("F3").Value = 0.3046272680949286
$ws.Range("G3").Value = 107.6780247903599

# Row 4
$ws.Range("D4").Value = "wRLS"
$ws.Range("E4").Value = 142.8100188340755
$ws.Range("F4").Value = 0.3026904009802125
$ws.Range("G4").Value = 107.0580792111746

# Row 5
$ws.Range("D5").Value = "wRLS"
$ws.Range("E5").Value = 142.5982545821705
$ws.Range("F5").Value = 0.3022415598775665
$ws.Range("G5").Value = 107.2101567648761

# Row 6
$ws.Range("D6").Value = "wRLS"
$ws.Range("E6").Value = 142.5264155072198
$ws.Range("F6").Value = 0.3020892946612993
$ws.Range("G6").Value = 107.2489954237841

# Row 7
$ws.Range("D7").Value = "wRLS"
$ws.Range("E7").Value = 142.3824198810079
$ws.Range("F7").Value = 0.3017840913275745
$ws.Range("G7").Value = 107.3954027425575

# Row 8
$ws.Range("D8").Value = "wRLS"
$ws.Range("E8").Value = 142.2962105181875
$ws.Range("F8").Value = 0.3016013678267068
$ws.Range("G8").Value = 107.3654978052027
